$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("630Adata")
$ws.Range("H9").Value = 9
